# Update the LR-pair TPM-derived statistics (Ntrk3-Ptprs) with new TPM-based
# recomputation. Only the receptor/edge expression & specificity columns
# (M..T) change; identifiers and counts (A..L) stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> FAPs)
$ws.Range("M2").Value = 2.425633666666667
$ws.Range("N2").Value = 7.276901000000001
$ws.Range("O2").Value = 0.0662600404061536
$ws.Range("P2").Value = 0.06626004040615362
$ws.Range("Q2").Value = 3.559635193813555
$ws.Range("R2").Value = 32.036716744322
$ws.Range("S2").Value = 0.01252814689587894
$ws.Range("T2").Value = 0.01252814689587894

# Row 3 (FAPs -> FAPs, different target)
$ws.Range("O3").Value = 0.4234968256437875
$ws.Range("P3").Value = 0.4234968256437876
$ws.Range("S3").Value = 0.08007285249272299
$ws.Range("T3").Value = 0.080072852492723

# Row 4 (FAPs -> MuSCs)
$ws.Range("M4").Value = 18.67887366666667
$ws.Range("N4").Value = 56.03662100000001
$ws.Range("O4").Value = 0.5102431339500588
$ws.Range("P4").Value = 0.5102431339500588
$ws.Range("Q4").Value = 27.41138408424023
$ws.Range("S4").Value = 0.09647444969179807
$ws.Range("T4").Value = 0.09647444969179807

# Row 5 (MuSCs -> FAPs/ECs)
$ws.Range("M5").Value = 2.425633666666667
$ws.Range("N5").Value = 7.276901000000001
$ws.Range("O5").Value = 0.0662600404061536
$ws.Range("P5").Value = 0.06626004040615362
$ws.Range("Q5").Value = 15.26689787077222
$ws.Range("R5").Value = 137.40208083695
$ws.Range("S5").Value = 0.05373189351027467
$ws.Range("T5").Value = 0.05373189351027467

# Row 6 (MuSCs -> FAPs)
$ws.Range("O6").Value = 0.4234968256437875
$ws.Range("P6").Value = 0.4234968256437876
$ws.Range("S6").Value = 0.3434239731510645
$ws.Range("T6").Value = 0.3434239731510645

# Row 7 (MuSCs -> MuSCs)
$ws.Range("M7").Value = 18.67887366666667
$ws.Range("N7").Value = 56.03662100000001
$ws.Range("O7").Value = 0.5102431339500588
$ws.Range("P7").Value = 0.5102431339500588
$ws.Range("S7").Value = 0.4137686842582608
$ws.Range("T7").Value = 0.4137686842582607
